$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 163, shifting the existing rows 163:254 down to 164:255.
$ws.Rows("163").Insert()

# Populate the newly inserted row 163 with the new weekly price-report entry.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R carry the same "template" values
# as the rest of the Femacal de La Calera / Ajo / Chino / Primera block, while
# D, J, K, L, M, P hold the genuinely new data for this row.
$ws.Range("A163").Value = 3
$ws.Range("B163").Value = "Femacal de La Calera"
$ws.Range("C163").Value = "Coquimbo"
$ws.Range("D163").Value = 44518
$ws.Range("E163").Value = 5
$ws.Range("F163").Value = 100112003
$ws.Range("G163").Value = "Ajo"
$ws.Range("H163").Value = "Chino"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 95
$ws.Range("K163").Value = 16000
$ws.Range("L163").Value = 16500
$ws.Range("M163").Value = 16263
$ws.Range("N163").Value = "$/caja 10 kilos"
$ws.Range("O163").Value = "China"
$ws.Range("P163").Value = 1626
$ws.Range("Q163").Value = 10
$ws.Range("R163").Value = "Hortaliza"
